# Update 13C-MFA files (run and result) for SC and IO under WT-batch and chemostats
#
# - Insert a new flux row "EX_glc__D_e.f" into FluxData (row 3), pushing the
#   existing rows down by one and updating the first few numeric values.
# - Make FluxData the active/selected sheet (instead of MSData), with the
#   view scrolled/zoomed/selected the way the saved workbook now shows it.

$wb = $excel.ActiveWorkbook

$msData = $wb.Worksheets.Item("MSData")
$fluxData = $wb.Worksheets.Item("FluxData")
$tracerData = $wb.Worksheets.Item("TracerData")

# ---------------------------------------------------------------------------
# FluxData: insert the new "EX_glc__D_e.f" row right after "BIOMASS.f" (row 2)
# ---------------------------------------------------------------------------
$fluxData.Rows.Item(3).Insert()

$fluxData.Range("A3").Value = "EX_glc__D_e.f"
$fluxData.Range("B3").Value = 2.1710936421265
$fluxData.Range("C3").Value = 0.889195720447356
$fluxData.Rows.Item(3).RowHeight = 13.8

# BIOMASS.f row keeps its flux value but gets a new (smaller) error value
$fluxData.Range("C2").Value = 0.0001

# The row that used to be "EX_c5sugal_e.f" (now row 4) gets new flux/error values
$fluxData.Range("B4").Value = 0.138888888888889
$fluxData.Range("C4").Value = 0.089785826002838

# ---------------------------------------------------------------------------
# View state: FluxData becomes the active sheet/tab, all sheets zoom to 110%,
# and selections collapse down to a single top-left cell.
# ---------------------------------------------------------------------------
[void]$msData.Activate()
[void]$msData.Range("A2").Select()
$excel.ActiveWindow.Zoom = 110

[void]$fluxData.Activate()
[void]$fluxData.Range("A1").Select()
$excel.ActiveWindow.Zoom = 110
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

[void]$tracerData.Activate()
[void]$tracerData.Range("A1").Select()
$excel.ActiveWindow.Zoom = 110

# FluxData is the sheet that should end up active/selected
[void]$fluxData.Activate()
[void]$fluxData.Range("A1").Select()
